# added condition for original date, done today date
#
# For rows 11-19 ("type N with ..." test rows), clear the leftover "No"
# placeholder answers in the Internal/Outside/Number-email columns (D/E/F)
# that shouldn't have been populated for those scenarios, and stamp the
# "Today date" column (I) with the run date (2022-12-03, serial 44898) so
# original-date can be conditioned on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Today date (column I) for rows 11-19 -------------------------------
# Register both the lowercase and uppercase custom date-only number format
# codes (mirrors the existing yyyy-mm-dd h:mm:ss / YYYY-MM-DD HH:MM:SS pair
# already in the workbook) by toggling the very first cell through both
# casings, then apply the final uppercase format to the rest directly so
# every cell in the column shares the same cell style.
$todayDate = 44898

$firstCell = $ws.Cells.Item(11, 9)
$firstCell.NumberFormat = "yyyy-mm-dd"
$firstCell.NumberFormat = "YYYY-MM-DD"
$firstCell.Value = $todayDate

for ($r = 12; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.NumberFormat = "YYYY-MM-DD"
    $cell.Value = $todayDate
}

# --- Clear stray "No" answers that no longer apply ----------------------
$ws.Range("D12").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("F19").Value = ""
